$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    Appears on the Overview sheet (zh-cn / de-de status columns E2, F2)
#    and on each per-locale sheet's "Status" column (C2).
# ---------------------------------------------------------------------------
$statusEdits = @(
    @{ Sheet = "Overview"; Cells = @("E2", "F2") },
    @{ Sheet = "zh-cn";    Cells = @("C2") },
    @{ Sheet = "de-de";    Cells = @("C2") }
)

foreach ($edit in $statusEdits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($addr in $edit.Cells) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Narrow the status columns (was ~17.22 chars, now ~13.41 chars) on the
#    Overview sheet (columns E & F) and on each locale sheet (column C).
#    12.5 is the ColumnWidth input that this engine's pixel-snapping rounds
#    to the stored width closest to the target 13.4101845877511.
# ---------------------------------------------------------------------------
$narrowWidth = 12.5

$widthEdits = @(
    @{ Sheet = "Overview"; Cols = @("E", "F") },
    @{ Sheet = "zh-cn";    Cols = @("C") },
    @{ Sheet = "de-de";    Cols = @("C") }
)

foreach ($edit in $widthEdits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($col in $edit.Cols) {
        $ws.Columns($col).ColumnWidth = $narrowWidth
    }
}
